$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.872.45"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.635.80"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5089"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2558"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06356"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07764"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "1.652.19"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.270"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5424"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "0.0₅7763"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "25.936.24"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "195.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.416"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.915"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.015"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.881"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1207"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.832"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04923"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.234"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.164"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.530"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.371"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8890"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.579"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Value = "1.127.70"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5407"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01547"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.004"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.543"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8127"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.536"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("D46").Value = "1.776.06"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4541"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.08%  "
